# Applies "Minor updates based on validation" changes:
#  - Adds "Blood Cell Count Ratio Measurements" category to the Blood Cell Count Ratio
#    Measurement BC's bc_categories list (row 48 in Biomedical Concepts, row 7 in BC Hierarchy)
#  - Adds "Body Measurements" category alongside "Vital Signs" for several Body Measurement BCs
#  - Adds "Blood Pressure" category alongside "Vital Signs" for several Blood Pressure BCs
#  - Fixes the "HematologyTests" typo / mis-categorization to "Chemistry Tests" for the
#    Prothrombin related BCs
#  - Keeps the "Categories" lookup sheet (unique, sorted list of all bc_categories values)
#    in sync with the above changes

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Biomedical Concepts" sheet - column F (bc_categories)
# ---------------------------------------------------------------------------
$wsBC = $wb.Worksheets.Item("Biomedical Concepts")

$bcUpdates = @{
    48  = "Laboratory Tests;Hematology Tests;Blood Cell Counts;Blood Cell Count Ratio Measurements"
    49  = "Vital Signs;Body Measurements"
    50  = "Vital Signs;Body Measurements"
    63  = "Vital Signs;Blood Pressure"
    64  = "Vital Signs;Blood Pressure"
    65  = "Vital Signs;Blood Pressure"
    66  = "Vital Signs;Blood Pressure"
    67  = "Vital Signs;Blood Pressure"
    68  = "Vital Signs;Blood Pressure"
    69  = "Vital Signs;Blood Pressure"
    160 = "Laboratory Tests;Chemistry Tests;Protein or Enzyme Type Measurements;Prothrombin Measurements;Coagulation Study"
    161 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    162 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    163 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    164 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    165 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    200 = "Vital Signs;Blood Pressure"
    201 = "Vital Signs;Blood Pressure"
    202 = "Vital Signs;Blood Pressure"
    203 = "Vital Signs;Blood Pressure"
    204 = "Vital Signs;Blood Pressure"
    205 = "Vital Signs;Blood Pressure"
    206 = "Vital Signs;Blood Pressure"
    222 = "Vital Signs;Body Measurements"
    223 = "Vital Signs;Body Measurements"
    224 = "Vital Signs;Body Measurements"
    225 = "Vital Signs;Body Measurements"
}

foreach ($row in $bcUpdates.Keys) {
    $wsBC.Cells.Item($row, 6).Value = $bcUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2. "BC Hierarchy" sheet - column E (bc_categories)
# ---------------------------------------------------------------------------
$wsHier = $wb.Worksheets.Item("BC Hierarchy")

$hierUpdates = @{
    7  = "Laboratory Tests;Hematology Tests;Blood Cell Counts;Blood Cell Count Ratio Measurements"
    8  = "Vital Signs;Body Measurements"
    13 = "Vital Signs;Blood Pressure"
    34 = "Laboratory Tests;Chemistry Tests;Protein or Enzyme Type Measurements;Prothrombin Measurements;Coagulation Study"
    35 = "Laboratory Tests;Chemistry Tests;Coagulation Study;Prothrombin Activity Measurements;Drug-Induced Liver Injury;DILI"
    42 = "Vital Signs;Blood Pressure"
    47 = "Vital Signs;Body Measurements"
}

foreach ($row in $hierUpdates.Keys) {
    $wsHier.Cells.Item($row, 5).Value = $hierUpdates[$row]
}

# ---------------------------------------------------------------------------
# 3. "Categories" sheet - column A, sorted unique list of all bc_categories
#    values. Three new categories are introduced above ("Blood Cell Count
#    Ratio Measurements", "Blood Pressure", "Body Measurements") and the
#    obsolete "HematologyTests" typo entry no longer appears anywhere, so the
#    sheet is rebuilt from the new set of unique values (53 - 1 + 3 = 55 rows).
# ---------------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item("Categories")

$finalCategories = @(
    "APACHE II",
    "APCH1",
    "Acute Physiology and Chronic Health Evaluation II Clinical Classification",
    "Adverse Events",
    "Allergen-induced Antibody Measurements",
    "Antibody Measurements",
    "Arterial Blood Gas Measurements",
    "Autoantibody Measurements",
    "Blood Cell Count Ratio Measurements",
    "Blood Cell Counts",
    "Blood Pressure",
    "Blood Protein Measurements",
    "Body Measurements",
    "COVID-19 Tests",
    "Chemistry Tests",
    "Choriogonadotropin Measurements",
    "Clinical Trial Attribute",
    "Clinical or Research Assessment Classification",
    "Coagulation Study",
    "Coombs Tests",
    "DILI",
    "Drug-Induced Liver Injury",
    "Events",
    "Factor III Measurements",
    "Gram Negative Bacteria Measurements",
    "HCG Measurements",
    "Hematology Tests",
    "Hormone Measurements",
    "Immunogenicity Specimen Assessments",
    "Immunoglobulin G Measurements",
    "Immunohematology Tests",
    "Immunology Tests",
    "Laboratory Tests",
    "Liver Function Tests",
    "MVAI",
    "MVAI1",
    "Medical Conditions",
    "Medical History Events",
    "Microbial-induced Antibody Measurement",
    "Microbiology Tests",
    "Modified Van Assche Index Clinical Classification",
    "Modified Van Assche Index Clinical Classification Question",
    "Presenting Conditions",
    "Protein or Enzyme Type Measurements",
    "Prothrombin Activity Measurements",
    "Prothrombin Measurements",
    "QRS",
    "QRS Instrument Questions",
    "Reported Events",
    "SARS-CoV-2 Tests",
    "Serology Tests",
    "Trial Summary",
    "Troponin Measurements",
    "Virology Tests",
    "Vital Signs"
)

$oldLastRow = 54   # previous UsedRange was A1:A54 (1 header + 53 values)
$newLastRow = 1 + $finalCategories.Count

# Extend formatting for any new rows beyond the previous extent by copying
# the format of the last pre-existing data row (A54) down before writing
# values, so the new rows pick up the same style used throughout column A.
if ($newLastRow -gt $oldLastRow) {
    $fmtSource = $wsCat.Range("A" + $oldLastRow)
    for ($r = $oldLastRow + 1; $r -le $newLastRow; $r++) {
        $fmtSource.Copy($wsCat.Range("A" + $r))
    }
}

for ($i = 0; $i -lt $finalCategories.Count; $i++) {
    $wsCat.Cells.Item(2 + $i, 1).Value = $finalCategories[$i]
}
